$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 865  # was 864
$ws.Range("F3").Value = 13838  # was 13826
$ws.Range("F4").Value = 13620  # was 13608
$ws.Range("F5").Value = 1055  # was 1054
$ws.Range("F7").Value = 43  # was 42
$ws.Range("F8").Value = 602  # was 601
$ws.Range("F9").Value = 82  # was 81
$ws.Range("F10").Value = 25  # was 24
$ws.Range("F12").Value = 769  # was 766
$ws.Range("F13").Value = 2151  # was 2150
$ws.Range("F14").Value = 112  # was 110
$ws.Range("F17").Value = 129  # was 126
$ws.Range("F19").Value = 535  # was 531
$ws.Range("F20").Value = 436  # was 434
$ws.Range("F21").Value = 412  # was 409
$ws.Range("F22").Value = 327  # was 325
$ws.Range("F23").Value = 270  # was 267
$ws.Range("F24").Value = 841  # was 838
$ws.Range("F25").Value = 97  # was 95
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 1536  # was 1529
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 113  # was 112
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 865  # was 864
$ws.Range("F4").Value = 13838  # was 13826
$ws.Range("F5").Value = 13620  # was 13608
$ws.Range("F6").Value = 1055  # was 1054
$ws.Range("F8").Value = 43  # was 42
$ws.Range("F9").Value = 602  # was 601
$ws.Range("F10").Value = 82  # was 81
$ws.Range("F11").Value = 25  # was 24
$ws.Range("F13").Value = 769  # was 766
$ws.Range("F16").Value = 2151  # was 2150
$ws.Range("F17").Value = 112  # was 110
$ws.Range("F20").Value = 129  # was 126
$ws.Range("F24").Value = 113  # was 112
$ws.Range("F25").Value = 113  # was 112
$ws.Range("F26").Value = 535  # was 531
$ws.Range("F27").Value = 436  # was 434
$ws.Range("F28").Value = 412  # was 409
$ws.Range("F29").Value = 327  # was 325
$ws.Range("F30").Value = 270  # was 267
$ws.Range("F31").Value = 841  # was 838
$ws.Range("F33").Value = 1536  # was 1529
$ws.Range("F37").Value = 97  # was 95
